$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing weekday hours for the week-of-43192 row (row 12).
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 2

# Update the active selection to match the author's final cursor position.
$ws.Range("M16").Select()
